# Cotações atualizadas - 2025-11-28
# Append the new day's quotes as row 84 (dates run A2:A83 = 45906..45988,
# this adds A84 = 45989, i.e. 2025-11-28).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Date column: reuse the same date/time number format used by the rest of
# column A (style index 2 in the original file -> "YYYY-MM-DD HH:MM:SS").
$ws.Range("A84").Value = 45989
$ws.Range("A84").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Quote columns are stored as plain text (comma decimal separator), same
# as every other data row in the sheet.
$ws.Range("B84").Value = "21,7883"
$ws.Range("C84").Value = "16,0515"
$ws.Range("D84").Value = "15,5122"
$ws.Range("E84").Value = "15,5122"
